$d = $word.ActiveDocument

$replacements = @(
    @{old = "681÷3="; new = "764÷6="},
    @{old = "281÷3="; new = "211÷7="},
    @{old = "910÷5="; new = "693÷9="},
    @{old = "766÷8="; new = "532÷8="},
    @{old = "761÷7="; new = "435÷7="},
    @{old = "827÷4="; new = "841÷7="},
    @{old = "630÷2="; new = "489÷7="},
    @{old = "353÷6="; new = "320÷8="},
    @{old = "656÷9="; new = "676÷8="},
    @{old = "507÷7="; new = "800÷2="},
    @{old = "119÷2="; new = "536÷9="},
    @{old = "230÷6="; new = "347÷8="},
    @{old = "161÷6="; new = "147÷6="},
    @{old = "692÷8="; new = "739÷4="},
    @{old = "544÷9="; new = "850÷7="},
    @{old = "594÷8="; new = "255÷6="},
    @{old = "621÷7="; new = "955÷5="},
    @{old = "875÷9="; new = "761÷2="},
    @{old = "419÷2="; new = "800÷4="},
    @{old = "872÷4="; new = "512÷6="},
    @{old = "700÷3="; new = "551÷9="},
    @{old = "325÷7="; new = "976÷9="},
    @{old = "362÷7="; new = "837÷6="},
    @{old = "374÷9="; new = "332÷7="},
    @{old = "438÷3="; new = "450÷9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
